# "copied over local copies of logs"
# Update the TASK SUMMARY SHEET hours and the ACTIVITY LOG SUMMARY SHEET
# totals to reflect the local (per-person) copies of the logs, and leave
# the workbook focused on the Activity Log Summary Sheet tab/cell the
# author was last looking at.

$wb = $excel.ActiveWorkbook

# --- TASK SUMMARY SHEET: "Hours Spent this Week" column (D) updates ---
$wsTask = $wb.Worksheets.Item("TASK SUMMARY SHEET")
$wsTask.Range("D3").Value = 6
$wsTask.Range("D4").Value = 2
$wsTask.Range("D5").Value = 7
$wsTask.Range("D6").Value = 5

# --- ACTIVITY LOG SUMMARY SHEET: Group/Individual work time updates ---
$wsSummary = $wb.Worksheets.Item("ACTIVITY LOG SUMMARY SHEET")
$wsSummary.Range("B4").Value = 3
$wsSummary.Range("C4").Value = ""
$wsSummary.Range("B5").Value = 12
$wsSummary.Range("B6").Value = 3
$wsSummary.Range("B7").Value = 2

# --- Restore each sheet's last selection ---
$wsLog = $wb.Worksheets.Item("ACTIVITY LOG SHEET ")
$wsLog.Activate()
$wsLog.Range("A22").Select()

$wsTask.Activate()
$wsTask.Range("D8").Select()

$wsSummary.Activate()
$wsSummary.Range("B9").Select()
